$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values scraped from the updated cryptos.xlsx export
$values = [ordered]@{
    'D2' = '29.439.28'
    'E2' = '  -0.12%  '
    'D3' = '1.852.33'
    'E3' = '  +0.04%  '
    'E4' = '  +0.14%  '
    'D5' = '240.89'
    'E5' = '  +0.04%  '
    'D6' = '0.6305'
    'E6' = '  -0.28%  '
    'E7' = '  +0.11%  '
    'D8' = '0.07677'
    'E8' = '  +1.40%  '
    'E9' = '  -0.75%  '
    'D10' = '24.60'
    'D11' = '0.07754'
    'E11' = '  +0.60%  '
    'D12' = '1.849.88'
    'E12' = '  -0.01%  '
    'D13' = '0.00001097'
    'E13' = '  +8.65%  '
    'E14' = '  +0.52%  '
    'D15' = '0.6811'
    'E15' = '  -0.68%  '
    'D16' = '83.62'
    'E16' = '  +0.35%  '
    'D17' = '2.105.91'
    'E17' = '  +0.66%  '
    'D18' = '6.154'
    'E18' = '  +0.06%  '
    'D19' = '29.468.12'
    'E19' = '  -0.14%  '
    'D20' = '229.71'
    'E20' = '  +0.28%  '
    'D21' = '12.47'
    'E21' = '  -0.35%  '
    'E22' = '  +0.14%  '
    'D23' = '7.454'
    'E23' = '  -1.13%  '
    'E24' = '  +0.10%  '
    'D25' = '156.88'
    'E25' = '  -0.03%  '
    'D26' = '0.1389'
    'E26' = '  -0.77%  '
    'D27' = '8.398'
    'E27' = '  +0.04%  '
    'D28' = '17.70'
    'E28' = '  +0.05%  '
    'E29' = '  +3.36%  '
    'E30' = '  -0.11%  '
    'D31' = '0.05723'
    'E31' = '  +0.57%  '
    'D32' = '4.135'
    'E32' = '  +0.05%  '
    'D33' = '4.050'
    'E33' = '  +0.34%  '
    'D34' = '1.852'
    'E34' = '  +0.27%  '
    'D35' = '1.163'
    'E35' = '  +0.28%  '
    'D36' = '0.7094'
    'E36' = '  -0.94%  '
    'D37' = '2.585'
    'E37' = '  -0.07%  '
    'E38' = '  +0.19%  '
    'E39' = '  -0.73%  '
    'D40' = '1.219.11'
    'E40' = '  -2.30%  '
    'D41' = '6.519'
    'E41' = '  +4.96%  '
    'D42' = '0.9086'
    'E42' = '  -0.36%  '
    'E43' = '  +0.12%  '
    'B44' = 'Quant'
    'C44' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D44' = '101.91'
    'E44' = '  +0.23%  '
    'B45' = 'Aave'
    'C45' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D45' = '66.48'
    'E45' = '  +0.51%  '
    'B46' = 'Aptos'
    'C46' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D46' = '7.133'
    'E46' = '  +0.72%  '
    'E47' = '  -0.41%  '
    'B48' = 'TheSandbox'
    'C48' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D48' = '0.4019'
    'E48' = '  -0.25%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D49' = '8.999'
    'E49' = '  -1.25%  '
    'B50' = 'RenderToken'
    'C50' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D50' = '1.684'
    'E50' = '  -0.64%  '
    'B51' = 'Algorand'
    'C51' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D51' = '0.1133'
    'E51' = '  +0.90%  '
}

# These Price-column cells hold numeric-looking text (e.g. "240.89").
# Force them to Text format first so Excel does not coerce the string
# into a floating point number (which would corrupt values such as
# "17.70" -> 17.7 or introduce binary floating point noise).
$textCells = @(
    'D5'
    'D6'
    'D8'
    'D10'
    'D11'
    'D13'
    'D15'
    'D16'
    'D18'
    'D20'
    'D21'
    'D23'
    'D25'
    'D26'
    'D27'
    'D28'
    'D31'
    'D32'
    'D33'
    'D34'
    'D35'
    'D36'
    'D37'
    'D41'
    'D42'
    'D44'
    'D45'
    'D46'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
